$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new rows 5-7, replicating the values of rows 2-4
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 10
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 20
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 30

# Copy the cell formatting (style) from A2:A4 onto the new A5:A7 cells
# so the new rows match the look of the existing data rows.
$ws.Range("A2:A4").Copy()
$ws.Range("A5:A7").PasteSpecial(-4122)
